$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings that must stay as text
# (matching formatting like "71.019.66", "7.50", "1.00", "8.80"), so we
# force the cell to Text format before assigning the value - otherwise
# Excel auto-converts them to numbers and drops significant trailing zeros.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.019.66"
$ws.Range("E2").Value = "  +5.75%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.636.72"
$ws.Range("E3").Value = "  +16.18%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.22"
$ws.Range("E5").Value = "  +2.88%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.52"
$ws.Range("E6").Value = "  +3.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.632.91"
$ws.Range("E7").Value = "  +16.07%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("E9").Value = "  +3.43%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("E10").Value = "  +6.69%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.62"
$ws.Range("E11").Value = "  +3.15%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.498"
$ws.Range("E12").Value = "  +4.28%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.88"
$ws.Range("E13").Value = "  +12.07%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000254"
$ws.Range("E14").Value = "  +4.63%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.239.08"
$ws.Range("E15").Value = "  +16.15%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "70.928.12"
$ws.Range("E16").Value = "  +5.74%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.640.99"
$ws.Range("E17").Value = "  +16.35%  "

# Row 18
$ws.Range("E18").Value = "  +0.87%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.50"
$ws.Range("E19").Value = "  +6.32%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.05"
$ws.Range("E20").Value = "  -0.20%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "513.50"
$ws.Range("E21").Value = "  +4.57%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.13"
$ws.Range("E22").Value = "  +16.73%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.742"
$ws.Range("E23").Value = "  +6.58%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.37"
$ws.Range("E24").Value = "  +4.03%  "

# Row 25
$ws.Range("E25").Value = "  +9.19%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.51"
$ws.Range("E26").Value = "  +5.29%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.01"
$ws.Range("E27").Value = "  +6.81%  "

# Row 28
$ws.Range("E28").Value = "  -0.09%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.54"
$ws.Range("E29").Value = "  +9.34%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.16"
$ws.Range("E30").Value = "  +1.29%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.77"
$ws.Range("E31").Value = "  +6.08%  "

# Row 32
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000110"
$ws.Range("E32").Value = "  +16.48%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.50"
$ws.Range("E33").Value = "  +11.60%  "

# Row 34
$ws.Range("E34").Value = "  +2.71%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.00%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.13"
$ws.Range("E36").Value = "  +8.12%  "

# Row 37
$ws.Range("E37").Value = "  +6.14%  "

# Row 38
$ws.Range("E38").Value = "  +10.79%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.17"
$ws.Range("E39").Value = "  +7.52%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.93"
$ws.Range("E40").Value = "  +2.90%  "

# Row 41
$ws.Range("E41").Value = "  +4.50%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.16"
$ws.Range("E42").Value = "  -6.57%  "

# Row 43
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.80"
$ws.Range("E43").Value = "  +5.15%  "

# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.121.66"
$ws.Range("E44").Value = "  +11.27%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "415.16"
$ws.Range("E45").Value = "  +10.32%  "

# Row 46
$ws.Range("E46").Value = "  +3.94%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.52"
$ws.Range("E47").Value = "  +14.10%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0368"
$ws.Range("E48").Value = "  +5.45%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.05"
$ws.Range("E49").Value = "  +2.17%  "

# Row 51
$ws.Range("E51").Value = "  +10.37%  "
